$d = $word.ActiveDocument

# Locate the paragraph that still contains the literal "TCNS" placeholder text
# (e.g. "Tribal Consultation Fees for TCNS _trans_ref_num_").
$target = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs($i)
    if ($p.Range.Text.Contains("Fees for TCNS")) {
        $target = $p
        break
    }
}

if ($target -ne $null) {
    $pStart = $target.Range.Start
    $pText = $target.Range.Text
    $idx = $pText.IndexOf("TCNS")

    if ($idx -ge 0) {
        # Move the (hidden) "_GoBack" bookmark to sit right before "TCNS ".
        # Adding it here first cleanly splits the existing run at that exact
        # spot without merging any of the surrounding, already-separate runs.
        $bmPos = $pStart + $idx
        $bmRange = $d.Range($bmPos, $bmPos)
        $d.Bookmarks.Add("_GoBack", $bmRange)

        # Remove the "TCNS " placeholder (the word plus its trailing space),
        # leaving "...Fees for " directly followed by the bookmark and the
        # "_trans_ref_num_" text.
        $delRange = $d.Range($bmPos, $bmPos + 5)
        $delRange.Text = ""
    }
}
